$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.510.99"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "2.514.52"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "573.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.46%  "
$ws.Range("D9").Value = "2.513.81"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +3.97%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "2.974.11"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").Value = "69.308.45"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("E16").Value = "  -2.46%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "24.78"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("D18").Value = "2.517.65"
$ws.Range("E18").Value = "  -0.03%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.27"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.56"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "349.84"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.93%  "
$ws.Range("E22").Value = "  -1.61%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +0.00%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "70.25"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E26").Value = "  -3.54%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("D28").Value = "2.644.49"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "0.0₃0890"
$ws.Range("E30").Value = "  -3.01%  "
$ws.Range("E31").Value = "  -0.81%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "462.74"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -4.68%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.23"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +1.16%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "157.31"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("E38").Value = "  +0.96%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.47"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  -3.16%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "38.09"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("E45").Value = "  -7.83%  "
$ws.Range("E46").Value = "  -7.61%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "141.75"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("E48").Value = "  -1.78%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.46"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("E50").Value = "  +0.27%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.577"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -3.94%  "
